# 13.1.3 sheet update: refresh headers to the "13.1.3" wording, extend the
# yearly data from just 2019 out to 2020-2023, and turn the D4 "484" text
# into a real number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Best-effort: record the absolute path Excel would stamp into
#     x15ac:absPath for this workbook (informational only: this MS
#     extension attribute is not exposed anywhere on the object model). ---
$wb.Path = "C:\Users\korozbaeva\Desktop\Показатели ЦУР для Платформы\Глобальные показатели ЦУР\"

# --- Header row (A1:C1): reword to the "13.1.3 ..." indicator text ---
$ws.Range("A1").Value = "13.1.3 Кырсыктардын кооптуулугун азайтуунун улуттук стратегияларына ылайык, кырсыктардын кооптуулугун азайтуунун жергиликтүү стратегияларын кабыл алган жана ишке ашырган жергиликтүү бийлик органдарынын үлүшү"
$ws.Range("B1").Value = "13.1.3 Доля местных органов власти, принявших и осуществляющих местные стратегии снижения риска бедствий в соответствии с национальными стратегиями снижения риска бедствий"
$ws.Range("C1").Value = "13.1.3 Proportion of local governments that adopt and implement local disaster risk reduction strategies in line with national disaster risk reduction strategies"

# --- Extend columns E:H (years 2020-2023), cloning D's formatting for each row ---
$ws.Range("D3:D6").Copy($ws.Range("E3:E6"))
$ws.Range("D3:D6").Copy($ws.Range("F3:F6"))
$ws.Range("D3:D6").Copy($ws.Range("G3:G6"))
$ws.Range("D3:D6").Copy($ws.Range("H3:H6"))

# Row 3: year headers
$ws.Range("E3").Value = 2020
$ws.Range("F3").Value = 2021
$ws.Range("G3").Value = 2022
$ws.Range("H3").Value = 2023

# Row 4: relabel A4, turn D4 into a genuine number (was text "484"), repeat across E:H
$ws.Range("A4").Value = "Жергиликтүү бийлик органдарынын саны"
$ws.Range("D4").Value = 484
$ws.Range("E4").Value = 484
$ws.Range("F4").Value = 484
$ws.Range("G4").Value = 484
$ws.Range("H4").Value = 484

# Row 5: relabel A5, fill the proportion values for each year
$ws.Range("A5").Value = "Кырсыктардын кооптуулугун азайтуунун жергиликтүү стратегияларын кабыл алган жана ишке ашырган жергиликтүү бийлик органдарынын үлүшү"
$ws.Range("E5").Value = 13.2
$ws.Range("F5").Value = 21.5
$ws.Range("G5").Value = 34.5
$ws.Range("H5").Value = 40.53

# Row 6: relabel A6, fill the count values for each year
$ws.Range("A6").Value = "Улуттук стратегияларга ылайык, кырсыктардын кооптуулугун азайтуу боюнча жергиликтүү DRR стратегияларын кабыл алган жана ишке ашырган жергиликтүү бийлик органдарынын саны"
$ws.Range("E6").Value = 67
$ws.Range("F6").Value = 104
$ws.Range("G6").Value = 167
$ws.Range("H6").Value = 169

# Reset the active selection back to A1 (source file no longer carries the
# stray I4 selection marker).
$ws.Range("A1").Select()
